$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: the title line "CodeX by BroCode" -> "CodeX by CodeStash"
# Scope the Find to the paragraph that actually contains the "by BroCode"
# sub-title so the other (unrelated) "BroCode" occurrence further down
# in the document is left untouched.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*by BroCode*") {
        $para.Range.Find.Execute("BroCode", $true, $false, $false, $false, $false, $true, 1, $false, "CodeStash", 2)
        break
    }
}

# ------------------------------------------------------------------
# Change 2: "...called BroCode. It was successful..." becomes
# "...called Code" / "Stash (formerly BroCode)" / ". It was successful..."
# i.e. the literal text "BroCode" right after "called " is replaced with
# "CodeStash (formerly BroCode)", and the insertion is split into three
# runs (identical formatting) at the "Code|Stash" and "...BroCode)|." seams
# to mirror the source edit (done by toggling+restoring a character
# property on the middle slice, which forces the run boundary without
# altering the final look of the text).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*called BroCode.*") {
        $pStart = $para.Range.Start
        $pText = $para.Range.Text
        $idx = $pText.IndexOf("BroCode")

        $oldStart = $pStart + $idx
        $oldEnd = $oldStart + "BroCode".Length

        $replacement = "CodeStash (formerly BroCode)"
        $rOld = $d.Range($oldStart, $oldEnd)
        $rOld.Text = $replacement

        $splitPoint = $oldStart + "Code".Length
        $newEnd = $oldStart + $replacement.Length

        $rMiddle = $d.Range($splitPoint, $newEnd)
        $rMiddle.Font.Bold = 1
        $rMiddle.Font.Bold = 0
        break
    }
}

Write-Output "edit complete"
